# Word COM-interop script implementing the VisionDocument.docx revision
# described in the commit "Final revision. Submitted to clients".
#
# Summary of changes:
#  - Introduce the "(MM)" abbreviation for "MATLAB Marina" / "Marina" and
#    use it consistently throughout the Introduction, Needs and
#    Requirements, and Product/Solution Overview sections.
#  - Fix the "Needs and Requirements" / "Accessibility" heading typos
#    (content unaffected aside from "Accessibillity" -> "Accessibility").
#  - Add a sentence to the "Prioritizing content" bullet describing what
#    non-critical content includes.
#  - Expand the Accessibility paragraph with the list of supported
#    browsers.
#  - Register the two new "ListLabel 1" / "ListLabel 2" character styles
#    used by the list bullets.

$d = $word.ActiveDocument

# Plain global find & replace (exact text, case sensitive, whole story).
function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# A Range covering just paragraph $n's text, excluding the trailing
# paragraph mark (so formatting no-ops below don't touch the pilcrow).
function Get-ParaRange($n) {
    $p = $d.Paragraphs.Item($n).Range
    return $d.Range($p.Start, $p.End - 1)
}

# Toggling a no-op Bold flip on a Range forces the interop layer to emit a
# clean, explicit (but empty) <w:rPr/> on the run(s) under that range and
# guarantees a run boundary exists exactly at the range's edges - this is
# how the runs get split/normalized to match the target layout without
# altering any visible formatting.
function Force-RunBreak($theRange) {
    $theRange.Font.Bold = $true
    $theRange.Font.Bold = $false
}

# Sequentially locate $text inside $cursor (a Range that tracks how far
# we've progressed through a paragraph) and normalize/split the run there.
# Advances $cursor.Start past the match so repeated substrings are found
# left-to-right in document order.
function Split-Next($cursor, $text) {
    $hit = $d.Range($cursor.Start, $cursor.End)
    $hit.Find.Execute($text, $true, $false, $false, $false, $false, `
                       $true, 0, $false, "", 0) | Out-Null
    if ($hit.Find.Found) {
        Force-RunBreak $hit
        $cursor.Start = $hit.End
    }
}

# =============================================================================
# Title (paragraph 1): split "MATLAB Marina Mobile Redesign" into two runs.
# =============================================================================
$cursor = Get-ParaRange 1
Split-Next $cursor "MATLAB Marina"
Force-RunBreak $cursor

# =============================================================================
# Introduction paragraph (paragraph 9)
# =============================================================================
Replace-Text "The Armstrong MATLAB Marina website" "The MATLAB Marina (MM) website"
Replace-Text "The Marina Mobile Redesign Project" "The (MM) Mobile Redesign Project"
Replace-Text "accessing Marina on handheld" "accessing (MM) on handheld"

$cursor = Get-ParaRange 9
Split-Next $cursor "The "
Split-Next $cursor "MATLAB Marina"
Split-Next $cursor " "
Split-Next $cursor "(MM) "
Split-Next $cursor "website serves as an online resource for students taking Engineering courses that require the use of the MATLAB programming language and associated software. It covers basic to advanced programming concepts for engineering students, offering video tutorials, example programs, and exercises. The "
Split-Next $cursor "(MM)"
Split-Next $cursor " Mobile Redesign Project aims to improve the experience of students accessing "
Split-Next $cursor "(MM)"
Force-RunBreak $cursor

# =============================================================================
# "Needs and Requirements" heading (paragraph 10) - runs merge, no text change
# =============================================================================
Replace-Text "Needs and Requirements" "Needs and Requirements"
$cursor = Get-ParaRange 10
Force-RunBreak $cursor

# =============================================================================
# Needs and Requirements body (paragraph 11)
# =============================================================================
Replace-Text "The current MATLAB Marina website" "The current (MM) website"
Replace-Text "The Marina website is no exception" "The (MM) website is no exception"
Replace-Text "accessing the current Marina website on mobile devices" "accessing the current (MM) website on mobile devices"

$cursor = Get-ParaRange 11
Split-Next $cursor "The current (MM) website is difficult to navigate on smartphones and tablets, as it was designed without these different media in mind. The use of such handheld mobile devices to access online content has eclipsed the use of large-screened devices like laptops and desktop computers, prompting the need for online content to be easily accessible from handheld devices. The "
Split-Next $cursor "(MM)"
Split-Next $cursor " website is no exception. Students accessing the current "
Split-Next $cursor "(MM)"
Force-RunBreak $cursor

# =============================================================================
# Product/Solution Overview body (paragraph 13)
# =============================================================================
Replace-Text "The Marina website will be redesigned" "The (MM) website will be redesigned"

$cursor = Get-ParaRange 13
Split-Next $cursor "The "
Split-Next $cursor "(MM)"
Force-RunBreak $cursor

# =============================================================================
# "New navigation" bullet (paragraph 15) - runs merge, no text change
# =============================================================================
Replace-Text "New navigation - The new website will feature a simplified navigation menu, which can be toggled in and out of view on mobile devices where screen width is restricted. " "New navigation - The new website will feature a simplified navigation menu, which can be toggled in and out of view on mobile devices where screen width is restricted. "
$cursor = Get-ParaRange 15
Force-RunBreak $cursor

# =============================================================================
# "Scaling" bullet (paragraph 16) - runs merge, no text change
# =============================================================================
Replace-Text "Scaling – The redesigned site will scale content to appear appropriately sized on smaller screens. " "Scaling – The redesigned site will scale content to appear appropriately sized on smaller screens. "
$cursor = Get-ParaRange 16
Force-RunBreak $cursor

# =============================================================================
# "Prioritizing content" bullet (paragraph 17): merge + append new sentence
# =============================================================================
Replace-Text "hidden or pushed to the bottom of the page." "hidden or pushed to the bottom of the page. Non-critical content would include the search function (would be integrated into the navigation), the relevant links sidebar, and the updates/blog sidebar. "

$cursor = Get-ParaRange 17
Split-Next $cursor "Prioritizing content - The redesigned layout of each page will ensure that critical content remains at the fore when screen size is limited, while non-critical content will be hidden or pushed to the bottom of the page. "
Split-Next $cursor "Non-critical content would include the search function (would be integrated into the navigation), the relevant links sidebar, and the updates/blog sidebar. "
Force-RunBreak $cursor

# =============================================================================
# Scope & Limitations body (paragraph 19) - runs merge, no text change
# =============================================================================
Replace-Text "The content of the website will remain as is. The new site will be maintainable by current methods, and will not introduce an interface for managing content. Steps will be taken to ensure that future content will be easily integrated into the redesigned site. " "The content of the website will remain as is. The new site will be maintainable by current methods, and will not introduce an interface for managing content. Steps will be taken to ensure that future content will be easily integrated into the redesigned site. "
$cursor = Get-ParaRange 19
Force-RunBreak $cursor

# =============================================================================
# Accessibility heading (paragraph 20): fix "Accessibillity" typo
# =============================================================================
Replace-Text "Accessibillity" "Accessibility"
$cursor = Get-ParaRange 20
Force-RunBreak $cursor

# =============================================================================
# Accessibility body (paragraph 21): append supported-browser list
# =============================================================================
Replace-Text "web browsers. " "web browsers (Internet Explorer 9 or better, Safari 5.1 or better, and current versions of Google Chrome, Mozilla Firefox, and Opera)."

$cursor = Get-ParaRange 21
Split-Next $cursor "The website will be accessible from all modern mobile and desktop/laptop web browsers "
Split-Next $cursor "(Internet Explorer 9 or better, Safari 5.1 or better, and current versions of Google Chrome, Mozilla Firefox, and Opera)."
Force-RunBreak $cursor

# =============================================================================
# Register the ListLabel character styles used by list bullets
# =============================================================================
$listLabel1 = $d.Styles.Add("ListLabel 1", 2)
$listLabel1.Font.NameBi = "Symbol"

$listLabel2 = $d.Styles.Add("ListLabel 2", 2)
$listLabel2.Font.NameBi = "OpenSymbol"
